# Auto-generated: apply market-data refresh edits per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1903.9
$ws.Range("I40").Value = 1903.9
$ws.Range("K40").Value = 1903.9
$ws.Range("M40").Value = -1728.9

$ws.Range("H58").Value = 754
$ws.Range("I58").Value = 90
$ws.Range("J58").Value = 1750
$ws.Range("K58").Value = 270
$ws.Range("L58").Value = 5250
$ws.Range("M58").Value = -120
$ws.Range("N58").Value = -5550

$ws.Range("H70").Value = 1523.125
$ws.Range("I70").Value = 1460
$ws.Range("K70").Value = 4380
$ws.Range("M70").Value = -4110

$ws.Range("H73").Value = 1523.125
$ws.Range("I73").Value = 1460
$ws.Range("K73").Value = 4380
$ws.Range("M73").Value = -3444

$ws.Range("H106").Value = 2578.3333
$ws.Range("I106").Value = 2171.2856
$ws.Range("K106").Value = 2171.2856
$ws.Range("M106").Value = -1540.2856

$ws.Range("H138").Value = 7873.3335
$ws.Range("J138").Value = 7873.3335
$ws.Range("L138").Value = 23620.0005
$ws.Range("N138").Value = -33900.00049999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H32").Value = 28594.5
$ws.Range("I32").Value = 26965.285
$ws.Range("K32").Value = 26965.285
$ws.Range("M32").Value = -26678.285

$ws.Range("H63").Value = 3875
$ws.Range("I63").Value = 3916.6667
$ws.Range("J63").Value = 3750
$ws.Range("K63").Value = 3916.6667
$ws.Range("L63").Value = 3750
$ws.Range("M63").Value = -3230.6667
$ws.Range("N63").Value = -5122

$ws.Range("H66").Value = 3875
$ws.Range("I66").Value = 3916.6667
$ws.Range("J66").Value = 3750
$ws.Range("K66").Value = 19583.3335
$ws.Range("L66").Value = 3750
$ws.Range("M66").Value = -16151.3335
$ws.Range("N66").Value = -25614

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H26").Value = 35499.5
$ws.Range("I26").Value = 35499.5
$ws.Range("K26").Value = 35499.5
$ws.Range("M26").Value = -35207.5

$ws.Range("H94").Value = 2467.1667
$ws.Range("I94").Value = 1198.25
$ws.Range("K94").Value = 1198.25
$ws.Range("M94").Value = -747.25

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 27.1
$ws.Range("J7").Value = 53.25
$ws.Range("L7").Value = 53.25
$ws.Range("N7").Value = -279.25

$ws.Range("H16").Value = 5722.5
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5722.5
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5722.5
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -6296.5

$ws.Range("H31").Value = 5561.6665
$ws.Range("J31").Value = 5319.6665
$ws.Range("L31").Value = 5319.6665
$ws.Range("N31").Value = -5909.6665

$ws.Range("H34").Value = 5561.6665
$ws.Range("J34").Value = 5319.6665
$ws.Range("L34").Value = 5319.6665
$ws.Range("N34").Value = -5723.6665

$ws.Range("H35").Value = 2378
$ws.Range("I35").Value = 722.4
$ws.Range("K35").Value = 722.4
$ws.Range("M35").Value = -428.4

$ws.Range("H59").Value = 35000
$ws.Range("I59").Value = 35000
$ws.Range("K59").Value = 35000
$ws.Range("M59").Value = -33855

$ws.Range("H103").Value = 29999.334
$ws.Range("I103").Value = 29999.334
$ws.Range("K103").Value = 29999.334
$ws.Range("M103").Value = -28827.334

$ws.Range("H106").Value = 21123.2
$ws.Range("J106").Value = 21123.2
$ws.Range("L106").Value = 21123.2
$ws.Range("N106").Value = -23647.2

$ws.Range("H113").Value = 5722.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5722.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5722.5
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10062.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 11.5
$ws.Range("I11").Value = 11.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 34.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 105.5
$ws.Range("N11").ClearContents()

$ws.Range("H130").Value = 3500
$ws.Range("I130").Value = 3500
$ws.Range("K130").Value = 10500
$ws.Range("M130").Value = -5480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H113").Value = 3266.6667
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 3400
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 3400
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -7740

$ws.Range("H122").Value = 1750.6
$ws.Range("I122").Value = 1750.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5251.799999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2801.799999999999
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 800
$ws.Range("J46").Value = 800
$ws.Range("L46").Value = 800
$ws.Range("N46").Value = -1176

$ws.Range("H58").Value = 34333.332
$ws.Range("I58").Value = 3250
$ws.Range("J58").Value = 49875
$ws.Range("K58").Value = 3250
$ws.Range("L58").Value = 49875
$ws.Range("M58").Value = -2990
$ws.Range("N58").Value = -50395

$ws.Range("H61").Value = 1500
$ws.Range("I61").Value = 1500
$ws.Range("K61").Value = 1500
$ws.Range("M61").Value = -1298

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 650

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H113").Value = 199.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 199.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 598.5
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -4938.5
